$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the 08:00-11:00 rule row from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Match the author's final cursor position / selection (E8) recorded in the saved file
$ws.Range("E8").Select()
